{"js": "// Replace the date line and every three-digit \u00f7 one-digit division problem\n// in the worksheet, in document order. Each old value in the table is\n// unique, and the diff maps each text run to exactly one new value in the\n// same order the runs appear in the document (title first, then each\n// non-empty table cell, row by row, left to right).\nconst replacements = [\n  [\"2025-07-17 Thursday\", \"2025-07-18 Friday\"],\n  [\"396\u00f74=99, 0\", \"191\u00f74=47, 3\"],\n  [\"944\u00f73=314, 2\", \"749\u00f78=93, 5\"],\n  [\"930\u00f79=103, 3\", \"226\u00f74=56, 2\"],\n  [\"147\u00f73=49, 0\", \"379\u00f73=126, 1\"],\n  [\"962\u00f76=160, 2\", \"118\u00f72=59, 0\"],\n  [\"671\u00f78=83, 7\", \"374\u00f75=74, 4\"],\n  [\"906\u00f73=302, 0\", \"797\u00f77=113, 6\"],\n  [\"555\u00f75=111, 0\", \"727\u00f76=121, 1\"],\n  [\"316\u00f72=158, 0\", \"385\u00f77=55, 0\"],\n  [\"556\u00f76=92, 4\", \"444\u00f78=55, 4\"],\n  [\"185\u00f76=30, 5\", \"847\u00f79=94, 1\"],\n  [\"504\u00f78=63, 0\", \"743\u00f73=247, 2\"],\n  [\"375\u00f72=187, 1\", \"718\u00f72=359, 0\"],\n  [\"422\u00f75=84, 2\", \"581\u00f79=64, 5\"],\n  [\"526\u00f77=75, 1\", \"345\u00f76=57, 3\"],\n  [\"254\u00f72=127, 0\", \"446\u00f79=49, 5\"],\n  [\"490\u00f74=122, 2\", \"985\u00f75=197, 0\"],\n  [\"514\u00f72=257, 0\", \"400\u00f76=66, 4\"],\n  [\"199\u00f72=99, 1\", \"427\u00f75=85, 2\"],\n  [\"200\u00f72=100, 0\", \"987\u00f77=141, 0\"],\n  [\"531\u00f75=106, 1\", \"519\u00f75=103, 4\"],\n  [\"152\u00f77=21, 5\", \"223\u00f73=74, 1\"],\n  [\"670\u00f78=83, 6\", \"755\u00f73=251, 2\"],\n  [\"545\u00f74=136, 1\", \"439\u00f75=87, 4\"],\n  [\"723\u00f72=361, 1\", \"688\u00f73=229, 1\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet nextReplacement = 0;\nfor (const paragraph of paragraphs.items) {\n  if (nextReplacement >= replacements.length) break;\n  const current = paragraph.text;\n  const [expectedOld, newText] = replacements[nextReplacement];\n  if (current === expectedOld) {\n    paragraph.insertText(newText, \"Replace\");\n    nextReplacement++;\n  }\n}\nawait context.sync();\n\nif (nextReplacement !== replacements.length) {\n  throw new Error(\n    \"Only replaced \" + nextReplacement + \" of \" + replacements.length + \" expected text runs\"\n  );\n}\n", "ps1": "# Replace the date line and every three-digit \u00f7 one-digit division problem\n# in the worksheet. Every \"old\" value that appears in the document is\n# unique, so a whole-document Find/Replace (wdReplaceAll, scoped to a single\n# hit each) for each pair is a safe, order-independent way to land the\n# diff's 26 text-run edits.\n$pairs = @(\n    @(\"2025-07-17 Thursday\", \"2025-07-18 Friday\"),\n    @(\"396\u00f74=99, 0\", \"191\u00f74=47, 3\"),\n    @(\"944\u00f73=314, 2\", \"749\u00f78=93, 5\"),\n    @(\"930\u00f79=103, 3\", \"226\u00f74=56, 2\"),\n    @(\"147\u00f73=49, 0\", \"379\u00f73=126, 1\"),\n    @(\"962\u00f76=160, 2\", \"118\u00f72=59, 0\"),\n    @(\"671\u00f78=83, 7\", \"374\u00f75=74, 4\"),\n    @(\"906\u00f73=302, 0\", \"797\u00f77=113, 6\"),\n    @(\"555\u00f75=111, 0\", \"727\u00f76=121, 1\"),\n    @(\"316\u00f72=158, 0\", \"385\u00f77=55, 0\"),\n    @(\"556\u00f76=92, 4\", \"444\u00f78=55, 4\"),\n    @(\"185\u00f76=30, 5\", \"847\u00f79=94, 1\"),\n    @(\"504\u00f78=63, 0\", \"743\u00f73=247, 2\"),\n    @(\"375\u00f72=187, 1\", \"718\u00f72=359, 0\"),\n    @(\"422\u00f75=84, 2\", \"581\u00f79=64, 5\"),\n    @(\"526\u00f77=75, 1\", \"345\u00f76=57, 3\"),\n    @(\"254\u00f72=127, 0\", \"446\u00f79=49, 5\"),\n    @(\"490\u00f74=122, 2\", \"985\u00f75=197, 0\"),\n    @(\"514\u00f72=257, 0\", \"400\u00f76=66, 4\"),\n    @(\"199\u00f72=99, 1\", \"427\u00f75=85, 2\"),\n    @(\"200\u00f72=100, 0\", \"987\u00f77=141, 0\"),\n    @(\"531\u00f75=106, 1\", \"519\u00f75=103, 4\"),\n    @(\"152\u00f77=21, 5\", \"223\u00f73=74, 1\"),\n    @(\"670\u00f78=83, 6\", \"755\u00f73=251, 2\"),\n    @(\"545\u00f74=136, 1\", \"439\u00f75=87, 4\"),\n    @(\"723\u00f72=361, 1\", \"688\u00f73=229, 1\")\n)\n\n$d = $word.ActiveDocument\n$missed = 0\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $found = $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $found) {\n        $missed = $missed + 1\n        Write-Output (\"MISSING: \" + $old)\n    }\n}\n\nif ($missed -gt 0) {\n    throw (\"Failed to find/replace \" + $missed + \" of \" + $pairs.Count + \" expected text runs\")\n}\n"}
